# Blind model names in the "Evaluations" sheet (column C) so reviewers
# cannot tell which model produced which translation.
#
# Mapping (per commit message):
#   claude-opus-4.5 -> Model A
#   gemini-3-pro    -> Model B
#   gpt-5.1         -> Model C
#   kimi-k2         -> Model D

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluations")

$blindingMap = @{
    "claude-opus-4.5" = "Model A"
    "gemini-3-pro"    = "Model B"
    "gpt-5.1"         = "Model C"
    "kimi-k2"         = "Model D"
}

# Data rows run from row 2 through row 49 (row 1 is the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 49 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $modelName = $cell.Value2
    if ($null -ne $modelName -and $blindingMap.ContainsKey([string]$modelName)) {
        $cell.Value2 = $blindingMap[[string]$modelName]
    }
}
